$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5446.077
$ws.Range("I40").Value = 1111
$ws.Range("J40").Value = 9161.857
$ws.Range("K40").Value = 1111
$ws.Range("L40").Value = 9161.857
$ws.Range("M40").Value = -936
$ws.Range("N40").Value = -9511.857

$ws.Range("H76").Value = 7568.2812
$ws.Range("I76").Value = 9171.166999999999
$ws.Range("K76").Value = 9171.166999999999
$ws.Range("M76").Value = -8856.166999999999

$ws.Range("H79").Value = 7568.2812
$ws.Range("I79").Value = 9171.166999999999
$ws.Range("K79").Value = 9171.166999999999
$ws.Range("M79").Value = -8079.166999999999

$ws.Range("H98").Value = 727.3226
$ws.Range("I98").Value = 502.24
$ws.Range("K98").Value = 502.24
$ws.Range("M98").Value = 995.76

$ws.Range("H113").Value = 4723.923
$ws.Range("I113").Value = 3352.75
$ws.Range("J113").Value = 5333.3335
$ws.Range("K113").Value = 3352.75
$ws.Range("L113").Value = 5333.3335
$ws.Range("M113").Value = -98.75
$ws.Range("N113").Value = -11841.3335

$ws.Range("H116").Value = 92106.75
$ws.Range("J116").Value = 4566.6665
$ws.Range("L116").Value = 4566.6665
$ws.Range("N116").Value = -11450.6665

$ws.Range("H122").Value = 727.3226
$ws.Range("I122").Value = 502.24
$ws.Range("K122").Value = 1506.72
$ws.Range("M122").Value = 943.28

$ws.Range("H132").Value = 2051.4583
$ws.Range("I132").Value = 1082.2881
$ws.Range("J132").Value = 6450
$ws.Range("K132").Value = 3246.8643
$ws.Range("L132").Value = 19350
$ws.Range("M132").Value = -716.8643000000002
$ws.Range("N132").Value = -24410

$ws.Range("H135").Value = 1555.1578
$ws.Range("I135").Value = 828.5714
$ws.Range("J135").Value = 1979
$ws.Range("K135").Value = 7457.1426
$ws.Range("L135").Value = 17811
$ws.Range("M135").Value = -4922.1426
$ws.Range("N135").Value = -22881

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2942.7222
$ws.Range("I61").Value = 2268.348
$ws.Range("J61").Value = 3443.0645
$ws.Range("K61").Value = 2268.348
$ws.Range("L61").Value = 3443.0645
$ws.Range("M61").Value = -2056.348
$ws.Range("N61").Value = -3867.0645

$ws.Range("H74").Value = 2214.5642
$ws.Range("I74").Value = 1324
$ws.Range("J74").Value = 3494.75
$ws.Range("K74").Value = 1324
$ws.Range("L74").Value = 3494.75
$ws.Range("M74").Value = -450
$ws.Range("N74").Value = -5242.75

$ws.Range("H77").Value = 2214.5642
$ws.Range("I77").Value = 1324
$ws.Range("J77").Value = 3494.75
$ws.Range("K77").Value = 6620
$ws.Range("L77").Value = 17473.75
$ws.Range("M77").Value = -2252
$ws.Range("N77").Value = -26209.75

$ws.Range("H110").Value = 1255.7587
$ws.Range("I110").Value = 1177.5
$ws.Range("J110").Value = 1501.7142
$ws.Range("K110").Value = 1177.5
$ws.Range("L110").Value = 1501.7142
$ws.Range("M110").Value = 867.5
$ws.Range("N110").Value = -5591.7142

$ws.Range("H136").Value = 2942.7222
$ws.Range("I136").Value = 2268.348
$ws.Range("J136").Value = 3443.0645
$ws.Range("K136").Value = 6805.044
$ws.Range("L136").Value = 10329.1935
$ws.Range("M136").Value = -4255.044
$ws.Range("N136").Value = -15429.1935

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = $null
$ws.Range("M22").Value = $null
$ws.Range("N22").Value = 0

$ws.Range("H46").Value = 5000
$ws.Range("J46").Value = 5000
$ws.Range("L46").Value = 5000
$ws.Range("N46").Value = -5596

$ws.Range("H80").Value = 337
$ws.Range("I80").Value = 128
$ws.Range("J80").Value = 432
$ws.Range("K80").Value = 128
$ws.Range("L80").Value = 432
$ws.Range("M80").Value = 870
$ws.Range("N80").Value = -2428

$ws.Range("H83").Value = 337
$ws.Range("I83").Value = 128
$ws.Range("J83").Value = 432
$ws.Range("K83").Value = 640
$ws.Range("L83").Value = 2160
$ws.Range("M83").Value = 4352
$ws.Range("N83").Value = -12144

$ws.Range("H105").Value = 1941.4193
$ws.Range("I105").Value = 1849.1111
$ws.Range("J105").Value = 2069.2307
$ws.Range("K105").Value = 1849.1111
$ws.Range("L105").Value = 2069.2307
$ws.Range("M105").Value = -102.1111000000001
$ws.Range("N105").Value = -5563.2307

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2586.9644
$ws.Range("I31").Value = 1407.3334
$ws.Range("J31").Value = 4279.478
$ws.Range("K31").Value = 1407.3334
$ws.Range("L31").Value = 4279.478
$ws.Range("M31").Value = -1112.3334
$ws.Range("N31").Value = -4869.478

$ws.Range("H34").Value = 2586.9644
$ws.Range("I34").Value = 1407.3334
$ws.Range("J34").Value = 4279.478
$ws.Range("K34").Value = 1407.3334
$ws.Range("L34").Value = 4279.478
$ws.Range("M34").Value = -1205.3334
$ws.Range("N34").Value = -4683.478

$ws.Range("H58").Value = 1421.6364
$ws.Range("I58").Value = 888.6539
$ws.Range("K58").Value = 888.6539
$ws.Range("M58").Value = -685.6539

$ws.Range("H136").Value = 1421.6364
$ws.Range("I136").Value = 888.6539
$ws.Range("K136").Value = 2665.9617
$ws.Range("M136").Value = -115.9616999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 4700
$ws.Range("I56").Value = 4700
$ws.Range("K56").Value = 4700
$ws.Range("M56").Value = -4170

$ws.Range("H116").Value = 4275.6113
$ws.Range("I116").Value = 708.7143
$ws.Range("K116").Value = 2126.1429
$ws.Range("M116").Value = 1315.8571

$ws.Range("H123").Value = 3158.3333
$ws.Range("J123").Value = 3158.3333
$ws.Range("L123").Value = 9474.999899999999
$ws.Range("N123").Value = -14374.9999

$ws.Range("H131").Value = 1012.4691
$ws.Range("I131").Value = 5240
$ws.Range("J131").Value = 905.44305
$ws.Range("K131").Value = 15720
$ws.Range("L131").Value = 2716.32915
$ws.Range("M131").Value = -10680
$ws.Range("N131").Value = -12796.32915

$ws.Range("H132").Value = 4842.125
$ws.Range("I132").Value = 949.8333
$ws.Range("J132").Value = 7177.5
$ws.Range("K132").Value = 8548.4997
$ws.Range("L132").Value = 64597.5
$ws.Range("M132").Value = -6018.4997
$ws.Range("N132").Value = -69657.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5405.816
$ws.Range("I70").Value = 3915.9246
$ws.Range("J70").Value = 8839.044
$ws.Range("K70").Value = 3915.9246
$ws.Range("L70").Value = 8839.044
$ws.Range("M70").Value = -3645.9246
$ws.Range("N70").Value = -9379.044

$ws.Range("H73").Value = 5405.816
$ws.Range("I73").Value = 3915.9246
$ws.Range("J73").Value = 8839.044
$ws.Range("K73").Value = 3915.9246
$ws.Range("L73").Value = 8839.044
$ws.Range("M73").Value = -2979.9246
$ws.Range("N73").Value = -10711.044

$ws.Range("H113").Value = 1266
$ws.Range("I113").Value = 979.9091
$ws.Range("J113").Value = 1790.5
$ws.Range("K113").Value = 979.9091
$ws.Range("L113").Value = 1790.5
$ws.Range("M113").Value = 1190.0909
$ws.Range("N113").Value = -6130.5
